$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('A1').Value = 'tights for womens'
$ws.Range('A2').Value = 'tights from victoria secret'
$ws.Range('A3').Value = 'tights large'
$ws.Range('A4').Value = 'tights medium'
$ws.Range('A5').Value = 'tights mens'
$ws.Range('A6').Value = 'tights opaque for women'
$ws.Range('A7').Value = 'tights over the knee effect'
$ws.Range('A8').Value = 'tights pack'
$ws.Range('A9').Value = 'tights pants for women'
$ws.Range('A10').Value = 'tights pastel'
$ws.Range('A11').Value = 'tights plus'
$ws.Range('A12').Value = 'tights plus size'
$ws.Range('A13').Value = 'tights pockets'
$ws.Range('A14').Value = 'tights pregnant women'
$ws.Range('A15').Value = 'tights sheer black'
$ws.Range('A16').Value = 'tights shiny'
$ws.Range('A17').Value = 'tights shorts for women'
$ws.Range('A18').Value = 'tights size 8'
$ws.Range('A19').Value = 'tights sport'
$ws.Range('A20').Value = 'tights sports women'
$ws.Range('A21').Value = 'tights tall'
$ws.Range('A22').Value = 'tights tall women'
$ws.Range('A23').Value = 'tights thick'
$ws.Range('A24').Value = 'tights thigh high'
$ws.Range('A25').Value = 'tights white stars running'
$ws.Range('A26').Value = 'tights with back seam'
$ws.Range('A27').Value = 'tights with braces'
$ws.Range('A28').Value = 'tights with holes'
$ws.Range('A29').Value = 'tights with knee'
$ws.Range('A30').Value = 'tights with pockets for women'
$ws.Range('A31').Value = 'tights with seam up back'
$ws.Range('A32').Value = 'tights with skirt'
$ws.Range('A33').Value = 'tights with wide waistband'
$ws.Range('A34').Value = 'tights women'
$ws.Range('A35').Value = 'tights women long'
$ws.Range('A36').Value = 'tights women running'
$ws.Range('A37').Value = 'tights women tall'
$ws.Range('A38').Value = 'tights womens'
$ws.Range('A39').Value = 'tights xl tall'
$ws.Range('A40').Value = 'time and true leggings'
$ws.Range('A41').Value = 'time limited safe'
$ws.Range('A42').Value = 'tissue puff balls'
$ws.Range('A43').Value = 'to rise in darkness'
$ws.Range('A44').Value = 'todays special sign'
$ws.Range('A45').Value = 'toghts women'
$ws.Range('A46').Value = 'tommie copper compression leggings women'
$ws.Range('A47').Value = 'tommie copper compression wear'
$ws.Range('A48').Value = 'tommie copper core band'
$ws.Range('A49').Value = 'tommie copper fit'
$ws.Range('A50').Value = 'tommie copper for women'
$ws.Range('A51').Value = 'tommie copper knee for women'
$ws.Range('A52').Value = 'tommie copper knee performance'
$ws.Range('A53').Value = 'tommie copper leg'
$ws.Range('A54').Value = 'tommie copper leggings'
$ws.Range('A55').Value = 'tommie copper leggings for women'
$ws.Range('A56').Value = 'tommie copper leggings with back support'
$ws.Range('A57').Value = 'tommie copper pants'
$ws.Range('A58').Value = 'tommie copper pants women'
$ws.Range('A59').Value = 'tommie copper performance'
$ws.Range('A60').Value = 'tommie copper performance knee'
$ws.Range('A61').Value = 'tommie copper shorts women'
$ws.Range('A62').Value = 'tommie copper size chart'
$ws.Range('A63').Value = 'tommie copper tights'
$ws.Range('A64').Value = 'tommie copper underwear'
$ws.Range('A65').Value = 'tommie copper underwear for men'
$ws.Range('A66').Value = 'tommie copper wear'
$ws.Range('A67').Value = 'tommie copper women'
$ws.Range('A68').Value = 'tommies underwear for women'
$ws.Range('A69').Value = 'tommy control leggings for women'
$ws.Range('A70').Value = 'tommy copper back support for women'
$ws.Range('A71').Value = 'tommy copper compression'
$ws.Range('A72').Value = 'tommy copper compression knee'
$ws.Range('A73').Value = 'tommy copper compression leggings'
$ws.Range('A74').Value = 'tommy copper compression pants'
$ws.Range('A75').Value = 'tommy copper compression tights'
$ws.Range('A76').Value = 'tommy copper compression wear'
$ws.Range('A77').Value = 'tommy copper core'
$ws.Range('A78').Value = 'tommy copper for women'
$ws.Range('A79').Value = 'tommy copper hamstring'
$ws.Range('A80').Value = 'tommy copper hip'
$ws.Range('A81').Value = 'tommy copper hip compression'
$ws.Range('A82').Value = 'tommy copper knee'
$ws.Range('A83').Value = 'tommy copper knee compression'
$ws.Range('A84').Value = 'tommy copper knee support'
$ws.Range('A85').Value = 'tommy copper leg'
$ws.Range('A86').Value = 'tommy copper leggings'
$ws.Range('A87').Value = 'tommy copper pants for women'
$ws.Range('A88').Value = 'tommy copper performance'
$ws.Range('A89').Value = 'tommy copper shorts'
$ws.Range('A90').Value = 'tommy copper shorts women'
$ws.Range('A91').Value = 'tommy copper tights for women'
$ws.Range('A92').Value = 'tommy copper women'
$ws.Range('A93').Value = 'tommy leggings'
$ws.Range('A94').Value = 'tommy leggings women'
$ws.Range('A95').Value = 'tommy pants for women'
$ws.Range('A96').Value = 'tommy tights for women'
$ws.Range('A97').Value = 'tone thighs'
$ws.Range('A98').Value = 'tone wear'
$ws.Range('A99').Value = 'top ball the outdoor pickleball'
$ws.Range('A100').Value = 'top bike race game'
